$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "QUANGHAI"
$ws.Range("C2").Value = "Quản Lí"
$ws.Range("E2").Value = "Phạm Hải Chấm Công"

$ws.Range("A3").Value = "HAIPHAM"
$ws.Range("E3").Value = "Quang Hải Chấm Công"

$ws.Range("E1:F1").ColumnWidth = 21.166666666666668

$ws.Range("E2").Select()
